# NYPD CompStat weekly report update
# - Bumps "Number" (week number) in header from 31 -> 32
# - Updates the reporting week dates
# - Updates the weekly crime statistics table (rows 16-31)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header rich text: "Volume 32   Number  31" -> "...Number  32"
# ---------------------------------------------------------------------------
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "32"

# ---------------------------------------------------------------------------
# Header rich text: reporting week dates
# "Report Covering the Week  7/28/2025  Through  8/3/2025"
#   -> "Report Covering the Week  8/4/2025  Through  8/10/2025"
# Edit right-to-left so earlier character offsets stay valid.
# ---------------------------------------------------------------------------
$weekCell = $ws.Range("C9")
$weekCell.Characters(47, 8).Text = "8/10/2025"
$weekCell.Characters(27, 9).Text = "8/4/2025"

# ---------------------------------------------------------------------------
# Helper functions
# ---------------------------------------------------------------------------
$fmtCount = "#,##0"
$fmtPct1  = '#,##0.0;"-"#,##0.0'

function Set-Count($addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-Pct($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Re-applies number format only when a cell is converting from the
# "N/A" text placeholder to a genuine numeric value, so the existing
# style (shared by every other numeric cell in that column) is reused
# instead of Excel minting a brand-new style entry.
function Set-CountFromText($addr, $val) {
    $ws.Range($addr).Value = $val
    $ws.Range($addr).NumberFormat = $fmtCount
}

function Set-PctFromText($addr, $val) {
    $ws.Range($addr).Value = $val
    $ws.Range($addr).NumberFormat = $fmtPct1
}

# Converts a numeric cell back into the "0" text placeholder used
# throughout the sheet for precincts with zero incidents, by copying
# an existing placeholder cell (keeps formatting + shared string type).
function Set-ZeroText($addr) {
    $ws.Range("C14").Copy($ws.Range($addr))
}

# ---------------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------------
Set-Count "C16" 2
Set-Count "D16" 4
Set-Pct   "E16" -50
Set-Count "F16" 11
Set-Count "G16" 12
Set-Pct   "H16" -8.333333333333
Set-Count "I16" 57
Set-Count "J16" 61
Set-Pct   "K16" -6.557377049180
Set-Pct   "L16" -1.724137931034
Set-Pct   "M16" 1.785714285714
Set-Pct   "N16" -84.718498659517

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
Set-Count "D17" 4
Set-Pct   "E17" -75
Set-Count "G17" 12
Set-Pct   "H17" -25
Set-Count "I17" 44
Set-Count "J17" 74
Set-Pct   "K17" -40.540540540540
Set-Pct   "L17" -35.294117647058
Set-Pct   "M17" 0
Set-Pct   "N17" -29.032258064516

# ---------------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------------
Set-Count "C18" 3
Set-Count "D18" 1
Set-Pct   "E18" 200
Set-Count "F18" 10
Set-Count "G18" 5
Set-Count "I18" 84
Set-Count "J18" 89
Set-Pct   "K18" -5.617977528089
Set-Pct   "L18" 15.068493150684
Set-Pct   "M18" 21.739130434782
Set-Pct   "N18" -84.090909090909

# ---------------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------------
Set-Count "C19" 18
Set-Count "D19" 20
Set-Pct   "E19" -10
Set-Count "F19" 62
Set-Count "G19" 56
Set-Pct   "H19" 10.714285714285
Set-Count "I19" 440
Set-Count "J19" 470
Set-Pct   "K19" -6.382978723404
Set-Pct   "L19" -0.900900900900
Set-Pct   "M19" 12.244897959183
Set-Pct   "N19" -62.457337883959

# ---------------------------------------------------------------------------
# Row 20 (D20/E20 flip from "N/A" text placeholders to real numbers)
# ---------------------------------------------------------------------------
Set-Count "C20" 1
Set-CountFromText "D20" 2
Set-PctFromText   "E20" -50
Set-Count "G20" 4
Set-Pct   "H20" 25
Set-Count "I20" 28
Set-Count "J20" 27
Set-Pct   "K20" 3.703703703703
Set-Pct   "L20" -59.420289855072
Set-Pct   "M20" 40
Set-Pct   "N20" -95.953757225433

# ---------------------------------------------------------------------------
# Row 21
# ---------------------------------------------------------------------------
Set-Count "C21" 25
Set-Count "D21" 31
Set-Pct   "E21" -19.354838709677
Set-Count "F21" 98
Set-Count "G21" 91
Set-Pct   "H21" 7.692307692307
Set-Count "I21" 659
Set-Count "J21" 728
Set-Pct   "K21" -9.478021978021
Set-Pct   "L21" -8.344923504867
Set-Pct   "M21" 11.884550084889
Set-Pct   "N21" -76.836555360281

# ---------------------------------------------------------------------------
# Row 22 (C22 flips from real number to "0" placeholder text; D22/E22
# flip from "N/A" text placeholders to real numbers)
# ---------------------------------------------------------------------------
Set-ZeroText "C22"
Set-CountFromText "D22" 1
Set-PctFromText   "E22" -100
Set-Count "G22" 2
Set-Pct   "H22" 0
Set-Count "J22" 16
Set-Pct   "K22" -18.75

# ---------------------------------------------------------------------------
# Row 23
# ---------------------------------------------------------------------------
Set-Count "G23" 1
Set-Pct   "H23" 0

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
Set-Count "C24" 45
Set-Count "D24" 21
Set-Pct   "E24" 114.285714285714
Set-Count "F24" 151
Set-Count "G24" 108
Set-Pct   "H24" 39.814814814814
Set-Count "I24" 880
Set-Count "J24" 783
Set-Pct   "K24" 12.388250319284
Set-Pct   "L24" 30.757800891530
Set-Pct   "M24" 41.025641025641

# ---------------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------------
Set-Count "C25" 33
Set-Count "D25" 18
Set-Pct   "E25" 83.333333333333
Set-Count "F25" 120
Set-Count "G25" 85
Set-Pct   "H25" 41.176470588235
Set-Count "I25" 664
Set-Count "J25" 620
Set-Pct   "K25" 7.096774193548
Set-Pct   "L25" 26.717557251908

# ---------------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------------
Set-Count "C26" 7
Set-Count "D26" 3
Set-Pct   "E26" 133.333333333333
Set-Count "F26" 23
Set-Count "G26" 18
Set-Pct   "H26" 27.777777777777
Set-Count "I26" 142
Set-Count "J26" 150
Set-Pct   "K26" -5.333333333333
Set-Pct   "L26" 0
Set-Pct   "M26" -17.441860465116

# ---------------------------------------------------------------------------
# Row 28 (C28 flips from real number to "0" placeholder text; D28/E28
# flip from "N/A" text placeholders to real numbers)
# ---------------------------------------------------------------------------
Set-ZeroText "C28"
Set-CountFromText "D28" 1
Set-PctFromText   "E28" -100
Set-Count "F28" 4
Set-Count "G28" 6
Set-Pct   "H28" -33.333333333333
Set-Count "J28" 18
Set-Pct   "K28" 44.444444444444
Set-Pct   "L28" 4

# ---------------------------------------------------------------------------
# Row 31
# ---------------------------------------------------------------------------
Set-Pct "L31" -28.571428571428
